$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in "Kraken2" -> "Kraken 2" header value (row 2) which is shared
# across the sheet; rows 3-13 reference it with a formula.
$ws.Range("Z2").Value = "Kraken 2"

# Rows 3 through 13 should reference the header cell via formula "=Z$2"
for ($row = 3; $row -le 13; $row++) {
    $ws.Range("Z$row").Formula = "=Z`$2"
}

# Update the active cell / selection to Z18
$ws.Range("Z18").Select()
